# ---------------------------------------------------------------------------
# Target change (per the supplied OOXML diff):
#   word/numbering.xml - four <w:abstractNum> definitions (abstractNumId
#   990, 99411, 991, 99414) each get their <w:nsid w:val="..."/> child
#   replaced with a new hex value:
#       990    : 2b0a00e4 -> d3bea296
#       99411  : 173fb073 -> 79dbc970
#       991    : e7f139ab -> b6f69cf0
#       99414  : 796088e7 -> 4d61c4b7
#   Nothing else differs anywhere in the package (no paragraph text,
#   formatting, list association/order, or any other attribute changes).
#
# Investigation performed before writing this script:
#   - `w:nsid` is the list "GUID" Word stamps on each abstractNum the first
#     time a list definition is minted; it is pure bookkeeping with no
#     visible or semantic effect, and — in real Word just as in this
#     iron_native host — it is NOT exposed anywhere in the Object Model:
#       * `$d.ListTemplates`, `List`, `ListFormat`, `ListLevel`, ... expose
#         no nsid/GUID member.
#       * `Range.WordOpenXML` / `Selection.WordOpenXML` *getters* do surface
#         the raw `word/numbering.xml` (nested in the pkg:package blob),
#         confirming the values are loaded/round-tripped, but the *setter*
#         is read-only in this host (raises "Document.WordOpenXML is a
#         read-only property ..." / silently no-ops via COM), and even on
#         real Word that setter only ever replaces the target Range's own
#         document-body content, never sibling package parts.
#       * `Range.InsertXML` likewise only rewrites the addressed range's
#         body content; any numbering/styles parts supplied alongside it
#         are used purely as read-only resolution context, not written
#         back (verified empirically: numbering.xml is byte-identical
#         after InsertXML calls that embed a modified numbering part).
#       * `Find.Execute` only ever matches visible document text, never
#         XML source/attributes (nsid values aren't part of any run's
#         text), so text find/replace cannot reach them.
#       * Reaching around the OM via the PowerShell host's own file I/O
#         (read/edit/rewrite the .docx zip's word/numbering.xml entry,
#         then `$word.Documents.Open()` it) is deliberately sealed off:
#         the active document's backing path is locked while open, and
#         `Documents.Open` only ever re-vends the already-loaded in-memory
#         document for a recognized/staged path (edits made to the bytes
#         on disk at that path are ignored), while any other path is
#         rejected outright as not among the staged input files.
#       * The only OM call that *can* mint/alter abstractNum/nsid data
#         (`ListFormat.ApplyListTemplateWithLevel`) can only append a
#         brand-new list template with engine-assigned content/id/nsid;
#         it cannot target/edit the existing abstractNumId 990/99411/
#         991/99414 entries in place, and using it would corrupt the
#         document with an unrelated extra list template and reassign a
#         paragraph's numbering — a much bigger (and wrong) change than
#         the one in the diff.
#   - The diff's four nsid swaps are the *only* differences in the whole
#     package; no document text changed (the word "cellulose" from the
#     commit message does not even occur in this document), so there is
#     no reachable, in-scope edit to perform through legitimate Word
#     automation — this matches the common case where a docx's numbering
#     part is incidentally re-minted by the tool that generated the file
#     (e.g. re-running a markdown->docx converter) rather than anything a
#     document author/automation script does.
#
# Conclusion: there is no Word COM-interop call — in this host or in real
# Word — capable of touching `w:nsid`. Making this script's output exactly
# match the diff is therefore out of reach of the supported object model,
# so this is intentionally left as a no-op that leaves every reachable
# (and everything else) part of the document untouched rather than
# faking an unrelated, corrupting edit.
$d = $word.ActiveDocument
